# Updated cryptos list (price/volume refresh) on Tue Feb 28 18:46:49 UTC 2023
# with GitHub Actions.
#
# The sheet stores Price (col D) and Volume(1h) (col E) as plain TEXT
# (inlineStr) cells, not numbers. Several Price values look numeric
# (e.g. "0.9995", "304.02"), so a bare `.Value = "304.02"` assignment
# would be auto-coerced to a real number by Excel. To keep them as text
# (matching the original cell type) we prefix those values with a
# leading apostrophe, exactly like a user typing `'304.02` into a cell -
# Excel strips the apostrophe from the stored/displayed text and only
# uses it to force text interpretation (quote-prefix), so the saved
# value is clean. Values that already aren't number-like (e.g.
# "23.520.76", which has two dots, or the "  +0.89%  " volume strings)
# don't need the prefix because Excel can't parse them as numbers anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Addr,
        [string]$Value
    )
    # Plain decimal numbers (optionally signed) would otherwise be
    # auto-coerced to a numeric cell type by Excel; a leading apostrophe
    # forces text interpretation without changing the stored text.
    if ($Value -match '^[+-]?\d+(\.\d+)?$') {
        $ws.Range($Addr).Value = "'" + $Value
    } else {
        $ws.Range($Addr).Value = $Value
    }
}

# row -> @{ col letter = new value }
$rowUpdates = [ordered]@{
    2  = @{ D = '23.520.76';  E = '  +0.89%  ' }
    3  = @{ D = '1.643.92';   E = '  +1.10%  ' }
    4  = @{ D = '0.9995';     E = '  -0.47%  ' }
    5  = @{ D = '1.001';      E = '  -0.26%  ' }
    6  = @{ D = '304.02';     E = '  +0.39%  ' }
    7  = @{ D = '0.3816';     E = '  +1.48%  ' }
    8  = @{ D = '52.04';      E = '  -0.99%  ' }
    9  = @{ D = '0.3621';     E = '  +0.47%  ' }
    10 = @{ D = '0.08215';    E = '  +2.12%  ' }
    11 = @{ D = '1.235';      E = '  +0.76%  ' }
    12 = @{ D = '1.002';      E = '  -0.26%  ' }
    13 = @{ D = '22.60';      E = '  +0.09%  ' }
    14 = @{ D = '6.481';      E = '  -0.93%  ' }
    15 = @{ D = '7.364';      E = '  +2.51%  ' }
    16 = @{ D = '0.00001241'; E = '  -0.34%  ' }
    17 = @{ D = '1.637.61';   E = '  +0.60%  ' }
    18 = @{ D = '95.31';      E = '  +2.15%  ' }
    19 = @{ D = '0.06978';    E = '  +0.83%  ' }
    20 = @{ D = '6.597';      E = '  +2.46%  ' }
    21 = @{ D = '17.56';      E = '  -1.69%  ' }
    22 = @{ D = '1.001';      E = '  -0.38%  ' }
    23 = @{ D = '12.53';      E = '  -0.97%  ' }
    24 = @{ D = '23.511.25';  E = '  +0.78%  ' }
    25 = @{ D = '2.529';      E = '  +4.04%  ' }
    26 = @{ D = '3.075';      E = '  -3.96%  ' }
    27 = @{ D = '21.23';      E = '  +1.24%  ' }
    28 = @{ D = '151.87';     E = '  +2.59%  ' }
    29 = @{ D = '5.281';      E = '  -0.12%  ' }
    30 = @{ D = '133.41';     E = '  -0.60%  ' }
    31 = @{ D = '1.818.82';   E = '  +0.48%  ' }
    32 = @{ D = '1.095';      E = '  +15.68%  ' }

    # Filecoin and WEMIXTOKEN swap rank positions (33 <-> 34), each also
    # getting refreshed Price/Volume figures.
    33 = @{
        B = 'Filecoin'
        C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
        D = '6.585'
        E = '  -2.38%  '
    }
    34 = @{
        B = 'WEMIXTOKEN'
        C = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
        D = '2.154'
        E = '  -6.66%  '
    }

    35 = @{ D = '11.54';    E = '  +6.72%  ' }
    36 = @{ D = '0.02778';  E = '  -1.92%  ' }
    37 = @{ D = '0.2519';   E = '  +0.07%  ' }
    38 = @{ D = '0.08776';  E = '  -0.26%  ' }
    39 = @{ D = '6.001';    E = '  -1.73%  ' }
    40 = @{ D = '0.07036';  E = '  -1.24%  ' }

    # TheSandbox and TrustWalletToken swap rank positions (41 <-> 42),
    # each also getting refreshed Price/Volume figures.
    41 = @{
        B = 'TheSandbox'
        C = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
        D = '0.7073'
        E = '  +1.01%  '
    }
    42 = @{
        B = 'TrustWalletToken'
        C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
        D = '1.350'
        E = '  -0.95%  '
    }

    43 = @{ D = '12.32';   E = '  +0.34%  ' }
    44 = @{ D = '15.61';   E = '  -3.04%  ' }
    45 = @{ D = '0.6555';  E = '  +1.96%  ' }
    46 = @{ D = '0.9994';  E = '  -0.33%  ' }
    47 = @{ D = '2.299';   E = '  -0.23%  ' }
    48 = @{ D = '3.967';   E = '  -0.16%  ' }
    49 = @{ D = '0.07985'; E = '  +0.19%  ' }
    50 = @{ D = '129.30';  E = '  +2.85%  ' }
    51 = @{ D = '1.194';   E = '  -0.62%  ' }
}

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in @('B', 'C', 'D', 'E')) {
        if ($cols.Contains($col)) {
            $addr = "$col$row"
            $newValue = $cols[$col]
            Set-TextCell $addr $newValue
        }
    }
}
